$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and row 3 in this sheet were two separate observation records that
# got swapped: row 2's record (Id, coordinates, comment) becomes what used
# to be row 3's, and vice versa. The Ålder-Stadium/Metod marker cells
# (K/L/N), which are empty-text placeholders present only on some rows,
# move from row 3 to row 2 as part of the swap.

# --- Id (column A) ---
$ws.Range("A2").Value = 131238351
$ws.Range("A3").Value = 131239702

# --- Ost / Nord (columns Q/R) ---
$ws.Range("Q2").Value = 464101
$ws.Range("R2").Value = 6758095

$ws.Range("Q3").Value = 464187
$ws.Range("R3").Value = 6758040

# --- Publik kommentar (column AC) ---
$ws.Range("AC2").Value = "Ringhack samt miljöbilder för området"
$ws.Range("AC3").Value = "2 bild, gran kåda"

# --- Empty-text placeholder cells K/L/N: move from row 3 to row 2 ---
foreach ($col in @("K2", "L2", "N2")) {
    $ws.Range($col).Value = "'"
    $ws.Range($col).Style = "Normal"
}

$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("N3").ClearContents()
